# The "empadronador" ranking table was refreshed with newer cumulative
# totals (total_registros) and re-sorted in descending order by that
# total. The set of people is unchanged - only the per-person totals
# increased and the row order was updated to reflect the new ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, already in the final desired (descending by total) order.
$data = @(
    @("FARFAN MONTOYA ROSSANA ISABEL", 63),
    @("DAVILA CORDOVA MARIBEL", 62),
    @("SANCHEZ ULLOA CESAR AUGUSTO", 57),
    @("RAMOS RAMOS HANDY JAIR", 57),
    @("GARCIA GUTIERREZ LUIS ARTURO", 55),
    @("FIORELA KEILY GUTIERREZ CRUZ", 55),
    @("CARBAJAL RAMOS JESUS MARINA", 53),
    @("ANGIE BELÉN RODRÍGUEZ ZAVALA", 53),
    @("GONZALES VALLE SEBASTIAN", 52),
    @("DELGADO DELGADO RONI", 51),
    @("OLIVA ALVA GOSSELYN NASSIRA", 50),
    @("BAZAN TEJADA JOSE VICENTE", 50),
    @("ARANEDA LOPEZ MARCO VIERI", 50),
    @("VERDE LIZARRAGA DEYSI EUFEMIA", 48),
    @("CASTILLO QUEZADA DIEGO ALONSO", 46),
    @("JOSSY IVANA SUÁREZ ZAVALETA", 44),
    @("GUERRA CALDERON ESTHEFANY NICOLLE", 43),
    @("DANY DARWIN VILLACORTA SAAVEDRA", 42),
    @("CARDENAS CAMPOJO MARY PAULA", 41),
    @("JAVE CHAVEZ ANGHELO MARTIN", 39),
    @("ALVITES CAMPOS SERGIO MARTIN", 20),
    @("ROCHA SIPIRAN JHORDAN ENRIQUE", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
